# Auto-generated edit script: replaces the data table in the 'Export' sheet
# with the refreshed Saldo export (re-sorted by balance, some accounts
# renumbered, a few accounts removed/added), per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 199,3
$data[0,0] = '004450724'
$data[0,1] = 'ASSAKO'
$data[0,2] = 3051174.5
$data[1,0] = '004472404'
$data[1,1] = 'DILSON'
$data[1,2] = 753585
$data[2,0] = '004451978'
$data[2,1] = 'ANTONIO'
$data[2,2] = 736201.19
$data[3,0] = '004332747'
$data[3,1] = 'LOHRAN'
$data[3,2] = 674202
$data[4,0] = '004488571'
$data[4,1] = 'CARLOS'
$data[4,2] = 585039.13
$data[5,0] = '004395314'
$data[5,1] = 'MARIA'
$data[5,2] = 459081.84
$data[6,0] = '004387250'
$data[6,1] = 'MONICA'
$data[6,2] = 212648.53
$data[7,0] = '004237325'
$data[7,1] = 'RICARDO'
$data[7,2] = 61381.08
$data[8,0] = '004222784'
$data[8,1] = 'RAFAEL'
$data[8,2] = 51352.38
$data[9,0] = '004644524'
$data[9,1] = 'PAULO'
$data[9,2] = 45632.63
$data[10,0] = '004334062'
$data[10,1] = 'MERG'
$data[10,2] = 26334.22
$data[11,0] = '004313254'
$data[11,1] = 'GUSTAVO'
$data[11,2] = 18656.83
$data[12,0] = '004364200'
$data[12,1] = 'BLOCO'
$data[12,2] = 14046.2
$data[13,0] = '004454365'
$data[13,1] = 'RAFAEL'
$data[13,2] = 13721.76
$data[14,0] = '004415557'
$data[14,1] = 'FILIPE'
$data[14,2] = 13634.89
$data[15,0] = '002277249'
$data[15,1] = 'DANILO'
$data[15,2] = 9735.91
$data[16,0] = '004361159'
$data[16,1] = 'HFR'
$data[16,2] = 6024.15
$data[17,0] = '004641487'
$data[17,1] = 'LAILA'
$data[17,2] = 4391.15
$data[18,0] = '004643737'
$data[18,1] = 'LARA'
$data[18,2] = 3345.55
$data[19,0] = '001761119'
$data[19,1] = 'BLUEMETRIX'
$data[19,2] = 992.14
$data[20,0] = '004563252'
$data[20,1] = 'FERNANDO'
$data[20,2] = 758.03
$data[21,0] = '004381180'
$data[21,1] = 'HFR'
$data[21,2] = 743.31
$data[22,0] = '004204500'
$data[22,1] = 'EDWARD'
$data[22,2] = 694.74
$data[23,0] = '004813088'
$data[23,1] = 'JULIANA'
$data[23,2] = 639.03
$data[24,0] = '005142592'
$data[24,1] = 'ALBERTO'
$data[24,2] = 551.25
$data[25,0] = '004804125'
$data[25,1] = 'EDUARDO'
$data[25,2] = 539.03
$data[26,0] = '004643746'
$data[26,1] = 'MARIO'
$data[26,2] = 506.63
$data[27,0] = '004119016'
$data[27,1] = 'HEMAT'
$data[27,2] = 399.92
$data[28,0] = '004216401'
$data[28,1] = 'SUELY'
$data[28,2] = 378.82
$data[29,0] = '004870976'
$data[29,1] = 'HFR'
$data[29,2] = 309.15
$data[30,0] = '005121919'
$data[30,1] = 'JORGE'
$data[30,2] = 297.95
$data[31,0] = '000626491'
$data[31,1] = 'FELIPE'
$data[31,2] = 280.01
$data[32,0] = '004436055'
$data[32,1] = 'MARCO'
$data[32,2] = 228.44
$data[33,0] = '004884046'
$data[33,1] = 'WILSON'
$data[33,2] = 219.38
$data[34,0] = '005022526'
$data[34,1] = 'ALEXANDRE'
$data[34,2] = 193.2
$data[35,0] = '004475395'
$data[35,1] = 'DAVID'
$data[35,2] = 185.02
$data[36,0] = '005020439'
$data[36,1] = 'BEATRIZ'
$data[36,2] = 179.56
$data[37,0] = '004754920'
$data[37,1] = 'LUIS'
$data[37,2] = 165.52
$data[38,0] = '004958578'
$data[38,1] = 'ASSAKO'
$data[38,2] = 145.19
$data[39,0] = '004556853'
$data[39,1] = 'MARCEL'
$data[39,2] = 139.41
$data[40,0] = '004956636'
$data[40,1] = 'BEATRIZ'
$data[40,2] = 119.14
$data[41,0] = '004382374'
$data[41,1] = 'THEOMAR'
$data[41,2] = 102.35
$data[42,0] = '004335144'
$data[42,1] = 'EDMUNDO'
$data[42,2] = 101.47
$data[43,0] = '002687737'
$data[43,1] = 'JOSE'
$data[43,2] = 100.02
$data[44,0] = '004908680'
$data[44,1] = 'ELENE'
$data[44,2] = 99.31
$data[45,0] = '004278033'
$data[45,1] = 'DAISY'
$data[45,2] = 97.63
$data[46,0] = '004339183'
$data[46,1] = 'JALISON'
$data[46,2] = 95.69
$data[47,0] = '004473942'
$data[47,1] = 'DAIANNE'
$data[47,2] = 95.09
$data[48,0] = '004477812'
$data[48,1] = 'DIEGO'
$data[48,2] = 95.02
$data[49,0] = '004334158'
$data[49,1] = 'LEONE'
$data[49,2] = 93.87
$data[50,0] = '004426743'
$data[50,1] = 'GABRIELLE'
$data[50,2] = 92.47
$data[51,0] = '004331477'
$data[51,1] = 'SUZY'
$data[51,2] = 91.78
$data[52,0] = '005206566'
$data[52,1] = 'LEVI'
$data[52,2] = 91
$data[53,0] = '004451996'
$data[53,1] = 'ADRIANO'
$data[53,2] = 90.54
$data[54,0] = '004212132'
$data[54,1] = 'JOAO'
$data[54,2] = 86.38
$data[55,0] = '004207278'
$data[55,1] = 'CESAR'
$data[55,2] = 85.29
$data[56,0] = '004360431'
$data[56,1] = 'CARLOS'
$data[56,2] = 85.25
$data[57,0] = '004277637'
$data[57,1] = 'LARA'
$data[57,2] = 84.69
$data[58,0] = '004862672'
$data[58,1] = 'RENATO'
$data[58,2] = 83.99
$data[59,0] = '004381194'
$data[59,1] = 'ALINNE'
$data[59,2] = 83.28
$data[60,0] = '004479734'
$data[60,1] = 'RODRIGO'
$data[60,2] = 83.06
$data[61,0] = '004472760'
$data[61,1] = 'SANDRA'
$data[61,2] = 82.78
$data[62,0] = '004452790'
$data[62,1] = 'GUSTAVO'
$data[62,2] = 81.8
$data[63,0] = '004332103'
$data[63,1] = 'JOSE'
$data[63,2] = 81.11
$data[64,0] = '004318604'
$data[64,1] = 'RENAN'
$data[64,2] = 80.51
$data[65,0] = '005173958'
$data[65,1] = 'VENIA'
$data[65,2] = 80.39
$data[66,0] = '004809902'
$data[66,1] = 'PEDRO'
$data[66,2] = 80.14
$data[67,0] = '004267976'
$data[67,1] = 'E3'
$data[67,2] = 79.84
$data[68,0] = '004350197'
$data[68,1] = 'GISELA'
$data[68,2] = 77.43
$data[69,0] = '004479287'
$data[69,1] = 'ANA'
$data[69,2] = 76.87
$data[70,0] = '004748761'
$data[70,1] = 'MARCELO'
$data[70,2] = 76.67
$data[71,0] = '004431591'
$data[71,1] = 'MARIO'
$data[71,2] = 73.59
$data[72,0] = '004713953'
$data[72,1] = 'ALESSANDRA'
$data[72,2] = 73.04
$data[73,0] = '004915243'
$data[73,1] = 'POLO'
$data[73,2] = 72.62
$data[74,0] = '005032151'
$data[74,1] = 'ANA'
$data[74,2] = 72.51
$data[75,0] = '004517506'
$data[75,1] = 'LUIZ'
$data[75,2] = 70.84
$data[76,0] = '004854514'
$data[76,1] = 'MARCIA'
$data[76,2] = 68.46
$data[77,0] = '004460491'
$data[77,1] = 'PEDRO'
$data[77,2] = 68.19
$data[78,0] = '004115403'
$data[78,1] = 'HEBERT'
$data[78,2] = 66.57
$data[79,0] = '004999434'
$data[79,1] = 'EDUARDO'
$data[79,2] = 65.85
$data[80,0] = '004855596'
$data[80,1] = 'MARIANA'
$data[80,2] = 64.36
$data[81,0] = '004254210'
$data[81,1] = 'MARCO'
$data[81,2] = 63.47
$data[82,0] = '004242237'
$data[82,1] = 'MARIAH'
$data[82,2] = 63
$data[83,0] = '004335251'
$data[83,1] = 'EDMUNDO'
$data[83,2] = 62.39
$data[84,0] = '005018038'
$data[84,1] = 'ELAINE'
$data[84,2] = 61.91
$data[85,0] = '004243043'
$data[85,1] = 'SUELI'
$data[85,2] = 59.23
$data[86,0] = '004340984'
$data[86,1] = 'RENATA'
$data[86,2] = 58.94
$data[87,0] = '005170415'
$data[87,1] = 'MONICA'
$data[87,2] = 58.93
$data[88,0] = '004452507'
$data[88,1] = 'DANIELA'
$data[88,2] = 57.99
$data[89,0] = '004581652'
$data[89,1] = 'CINCO'
$data[89,2] = 57.64
$data[90,0] = '004329229'
$data[90,1] = 'GABRIEL'
$data[90,2] = 56.99
$data[91,0] = '004215217'
$data[91,1] = 'CAROLINA'
$data[91,2] = 55.66
$data[92,0] = '001759765'
$data[92,1] = 'NATAL'
$data[92,2] = 54.77
$data[93,0] = '004321092'
$data[93,1] = 'DANIEL'
$data[93,2] = 53.54
$data[94,0] = '004866753'
$data[94,1] = 'GENESI'
$data[94,2] = 53.52
$data[95,0] = '004259650'
$data[95,1] = 'BENTO'
$data[95,2] = 51.65
$data[96,0] = '005009922'
$data[96,1] = 'ANA'
$data[96,2] = 51.64
$data[97,0] = '004998717'
$data[97,1] = 'GIOVANE'
$data[97,2] = 51.01
$data[98,0] = '004208447'
$data[98,1] = 'LEILA'
$data[98,2] = 50
$data[99,0] = '004392159'
$data[99,1] = 'RODRIGO'
$data[99,2] = 48.16
$data[100,0] = '003115072'
$data[100,1] = 'VICTOR'
$data[100,2] = 48.04
$data[101,0] = '005216881'
$data[101,1] = 'RENAN'
$data[101,2] = 46.76
$data[102,0] = '001294033'
$data[102,1] = 'VIVIANE'
$data[102,2] = 45.35
$data[103,0] = '004335031'
$data[103,1] = 'EDMUNDO'
$data[103,2] = 45.3
$data[104,0] = '004491730'
$data[104,1] = 'DENISE'
$data[104,2] = 45.11
$data[105,0] = '001731007'
$data[105,1] = 'GUILHERME'
$data[105,2] = 44.59
$data[106,0] = '004259659'
$data[106,1] = 'BENTO'
$data[106,2] = 43.81
$data[107,0] = '004360430'
$data[107,1] = 'VIOMAR'
$data[107,2] = 43.24
$data[108,0] = '004805133'
$data[108,1] = 'PATRICIA'
$data[108,2] = 41.48
$data[109,0] = '002697806'
$data[109,1] = 'CLAUDIA'
$data[109,2] = 40.23
$data[110,0] = '004238164'
$data[110,1] = 'DANIELA'
$data[110,2] = 38.3
$data[111,0] = '005103059'
$data[111,1] = 'WALQUIRIA'
$data[111,2] = 38.22
$data[112,0] = '004224815'
$data[112,1] = 'GUILHERME'
$data[112,2] = 36.48
$data[113,0] = '005000656'
$data[113,1] = 'LUCIA'
$data[113,2] = 35.88
$data[114,0] = '004643880'
$data[114,1] = 'GABRIEL'
$data[114,2] = 34.75
$data[115,0] = '004752615'
$data[115,1] = 'LUZIMAR'
$data[115,2] = 33.83
$data[116,0] = '005028018'
$data[116,1] = 'ALEXANDRE'
$data[116,2] = 33.05
$data[117,0] = '004397124'
$data[117,1] = 'MURYLO'
$data[117,2] = 32.73
$data[118,0] = '005133039'
$data[118,1] = 'PAULO'
$data[118,2] = 31.7
$data[119,0] = '004481463'
$data[119,1] = 'MARA'
$data[119,2] = 29.77
$data[120,0] = '004265173'
$data[120,1] = 'JULIA'
$data[120,2] = 28.88
$data[121,0] = '004486497'
$data[121,1] = 'ELENA'
$data[121,2] = 28.63
$data[122,0] = '002064834'
$data[122,1] = 'RAFAELA'
$data[122,2] = 26.6
$data[123,0] = '002738211'
$data[123,1] = 'MARGARETH'
$data[123,2] = 26.3
$data[124,0] = '001000882'
$data[124,1] = 'AYRTON'
$data[124,2] = 25.52
$data[125,0] = '004643153'
$data[125,1] = 'CARLA'
$data[125,2] = 25.37
$data[126,0] = '004504449'
$data[126,1] = 'KELMA'
$data[126,2] = 24.88
$data[127,0] = '004404724'
$data[127,1] = 'LEANDRO'
$data[127,2] = 24.14
$data[128,0] = '004755204'
$data[128,1] = 'FABIANA'
$data[128,2] = 22.07
$data[129,0] = '005046919'
$data[129,1] = 'MARIANA'
$data[129,2] = 22
$data[130,0] = '005245032'
$data[130,1] = 'ROSA'
$data[130,2] = 21.25
$data[131,0] = '004332207'
$data[131,1] = 'IRACY'
$data[131,2] = 20.83
$data[132,0] = '004214604'
$data[132,1] = 'MARIA'
$data[132,2] = 20.72
$data[133,0] = '004493324'
$data[133,1] = 'DANIEL'
$data[133,2] = 20.01
$data[134,0] = '004228456'
$data[134,1] = 'FLASH'
$data[134,2] = 19.56
$data[135,0] = '004497825'
$data[135,1] = 'PRISCILLA'
$data[135,2] = 19.49
$data[136,0] = '004204255'
$data[136,1] = 'AMADO'
$data[136,2] = 18.77
$data[137,0] = '004368994'
$data[137,1] = 'CRISTINA'
$data[137,2] = 18.56
$data[138,0] = '004399832'
$data[138,1] = 'EULER'
$data[138,2] = 17.43
$data[139,0] = '005274028'
$data[139,1] = 'RAFAEL'
$data[139,2] = 16.72
$data[140,0] = '005143579'
$data[140,1] = 'GABRIEL'
$data[140,2] = 16.18
$data[141,0] = '005169333'
$data[141,1] = 'EDUARDO'
$data[141,2] = 16.12
$data[142,0] = '004268684'
$data[142,1] = 'PATRICIA'
$data[142,2] = 15.41
$data[143,0] = '004213943'
$data[143,1] = 'ELISA'
$data[143,2] = 15.13
$data[144,0] = '004422594'
$data[144,1] = 'WANDIR'
$data[144,2] = 14.67
$data[145,0] = '000827730'
$data[145,1] = 'LUCIANA'
$data[145,2] = 13.29
$data[146,0] = '004752461'
$data[146,1] = 'SERGIO'
$data[146,2] = 10.77
$data[147,0] = '001719494'
$data[147,1] = 'LUIS'
$data[147,2] = 10.24
$data[148,0] = '004458604'
$data[148,1] = 'FABIOLA'
$data[148,2] = 9.91
$data[149,0] = '004216298'
$data[149,1] = 'FLORDELIZ'
$data[149,2] = 9.74
$data[150,0] = '004646727'
$data[150,1] = 'RENATA'
$data[150,2] = 9.1
$data[151,0] = '004921978'
$data[151,1] = 'ELAINE'
$data[151,2] = 8.08
$data[152,0] = '004381415'
$data[152,1] = 'JOAO'
$data[152,2] = 7.95
$data[153,0] = '004693631'
$data[153,1] = 'NELY'
$data[153,2] = 7.36
$data[154,0] = '004470679'
$data[154,1] = 'RODOLFO'
$data[154,2] = 7.35
$data[155,0] = '005228239'
$data[155,1] = 'DEBORA'
$data[155,2] = 7.11
$data[156,0] = '004530494'
$data[156,1] = 'ROSANGELA'
$data[156,2] = 6.94
$data[157,0] = '004854496'
$data[157,1] = 'JOSE'
$data[157,2] = 6.64
$data[158,0] = '004448501'
$data[158,1] = 'JOAO'
$data[158,2] = 5.55
$data[159,0] = '005142624'
$data[159,1] = 'RODRIGO'
$data[159,2] = 4.75
$data[160,0] = '004994036'
$data[160,1] = 'BALTASAR'
$data[160,2] = 4.67
$data[161,0] = '004207658'
$data[161,1] = 'ROBERTO'
$data[161,2] = 4.54
$data[162,0] = '004848927'
$data[162,1] = 'ULDARICO'
$data[162,2] = 3.62
$data[163,0] = '005142661'
$data[163,1] = 'SABRINA'
$data[163,2] = 3.6
$data[164,0] = '004425261'
$data[164,1] = 'THAYSA'
$data[164,2] = 3.59
$data[165,0] = '004945161'
$data[165,1] = 'SONIA'
$data[165,2] = 3
$data[166,0] = '004240292'
$data[166,1] = 'MARCO'
$data[166,2] = 2.66
$data[167,0] = '004451652'
$data[167,1] = 'MATEUS'
$data[167,2] = 2.64
$data[168,0] = '004886366'
$data[168,1] = 'RENATO'
$data[168,2] = 1.57
$data[169,0] = '004520100'
$data[169,1] = 'ALEXANDRE'
$data[169,2] = 1.46
$data[170,0] = '004264780'
$data[170,1] = 'MARCELO'
$data[170,2] = 1.42
$data[171,0] = '004308815'
$data[171,1] = 'ZELI'
$data[171,2] = 1.25
$data[172,0] = '002694089'
$data[172,1] = 'VITOR'
$data[172,2] = 1.12
$data[173,0] = '004459875'
$data[173,1] = 'HELVECIO'
$data[173,2] = 1.08
$data[174,0] = '005019925'
$data[174,1] = 'ALEXANDRE'
$data[174,2] = 1.06
$data[175,0] = '004452597'
$data[175,1] = 'LARA'
$data[175,2] = 0.89
$data[176,0] = '004214460'
$data[176,1] = 'MARIA'
$data[176,2] = 0.79
$data[177,0] = '004223502'
$data[177,1] = 'BRUNA'
$data[177,2] = 0.78
$data[178,0] = '004380749'
$data[178,1] = 'ELEUSE'
$data[178,2] = 0.6
$data[179,0] = '005165116'
$data[179,1] = 'ANA'
$data[179,2] = 0.51
$data[180,0] = '004862677'
$data[180,1] = 'RENATO'
$data[180,2] = 0.43
$data[181,0] = '004453302'
$data[181,1] = 'ISABELLA'
$data[181,2] = 0.39
$data[182,0] = '004924605'
$data[182,1] = 'ESTER'
$data[182,2] = 0.3
$data[183,0] = '004587511'
$data[183,1] = 'CARLOS'
$data[183,2] = 0.24
$data[184,0] = '004239387'
$data[184,1] = 'LUIZ'
$data[184,2] = 0.22
$data[185,0] = '004806286'
$data[185,1] = 'VERA'
$data[185,2] = 0.19
$data[186,0] = '004371857'
$data[186,1] = 'NAZARETH'
$data[186,2] = 0.18
$data[187,0] = '004357159'
$data[187,1] = 'JOAO'
$data[187,2] = 0.15
$data[188,0] = '004320840'
$data[188,1] = 'NATALIA'
$data[188,2] = 0.14
$data[189,0] = '004466350'
$data[189,1] = 'RAQUEL'
$data[189,2] = 0.11
$data[190,0] = '004806244'
$data[190,1] = 'CARLA'
$data[190,2] = 0.1
$data[191,0] = '005047946'
$data[191,1] = 'GABRIEL'
$data[191,2] = 0.09
$data[192,0] = '004213929'
$data[192,1] = 'RODOLFO'
$data[192,2] = 0.08
$data[193,0] = '004589311'
$data[193,1] = 'CLARICE'
$data[193,2] = 0.06
$data[194,0] = '004473718'
$data[194,1] = 'LUCAS'
$data[194,2] = 0.02
$data[195,0] = '004850070'
$data[195,1] = 'RENATO'
$data[195,2] = 0.02
$data[196,0] = '002878817'
$data[196,1] = 'GUILHERME'
$data[196,2] = 0.01
$data[197,0] = '004400000'
$data[197,1] = 'VILMA'
$data[197,2] = 0.01
$data[198,0] = '004999410'
$data[198,1] = 'SONIA'
$data[198,2] = -182.89

# Account numbers must stay text (preserve leading zeros) while the
# balance column must remain numeric, so format column A as Text
# before writing values (formatting the whole 3-col range would also
# force column C to text, which we don't want).
$ws.Range("A2:A200").NumberFormat = "@"
$ws.Range("A2:C200").Value = $data

# Clear the now-unused tail rows left over from the larger previous table
# (the refreshed export has 2 fewer data rows than before); this also
# clears the blank separator row and the old footer row so they can be
# rewritten at their new (shifted-up) position below.
$ws.Range("A201:C204").ClearContents()

# Blank separator row stays empty; write the filter-description footer
# row (column A only) right after it.
$ws.Range("A202").Value = 'Filtros aplicados:
DataFim é (Em branco)
nr_saldo_disponivel não é 0
Posição é Posição D-1
DataFim é (Em branco)
CARTEIRA não está em branco
NR_CONTA não está em branco
TIPO_LANCAMENTO não é ED, ET ou Liquidação Doador'
